$wb = $excel.ActiveWorkbook

# --- Sheet 1: "VENTAS POR GRUPO" (dimension A1:R6 -> A1:R7) ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert a new row above the current totals row (row 6); the old row 6
# (the "x de 4" summary row) shifts down to row 7, carrying its formatting.
$ws1.Rows.Item(6).Insert()

# Fill the newly inserted row 6 with the new advisor's data.
$ws1.Cells.Item(6, 1).Value = "OFICINA-CATAECSA"
$ws1.Cells.Item(6, 2).Value = "MOROCHO PLAZA SHIRLEY AURELIA"
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item(6, $col).Value = 0
}
# Match the numeric-cell formatting used by the other data rows (row 2..5).
$ws1.Range($ws1.Cells.Item(2, 3), $ws1.Cells.Item(2, 18)).Copy()
$ws1.Range($ws1.Cells.Item(6, 3), $ws1.Cells.Item(6, 18)).PasteSpecial(-4122)

# Update the summary row (now row 7): "x de 4" -> "x de 5".
for ($col = 3; $col -le 18; $col++) {
    $cell = $ws1.Cells.Item(7, $col)
    $cell.Value = ($cell.Value2 -replace "de 4", "de 5")
}

# --- Sheet 2: "VENTA MENSUAL" (dimension A1:G6 -> A1:G7) ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(6).Insert()

$ws2.Cells.Item(6, 1).Value = "OFICINA-CATAECSA"
$ws2.Cells.Item(6, 2).Value = "MOROCHO PLAZA SHIRLEY AURELIA"
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item(6, $col).Value = 0
}
$ws2.Range($ws2.Cells.Item(2, 3), $ws2.Cells.Item(2, 7)).Copy()
$ws2.Range($ws2.Cells.Item(6, 3), $ws2.Cells.Item(6, 7)).PasteSpecial(-4122)

$excel.CutCopyMode = 0
